$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").NumberFormat = "@"
$ws.Range("A1").Value = "3273114"
$ws.Range("A1").Font.Bold = $true
$ws.Range("A1").WrapText = $true
$ws.Range("A1").NumberFormat = "General"

$ws.Range("B1").ClearContents()
$ws.Range("B1").Font.Bold = $true

$ws.Range("A2").Value = "6SL32105BE211UV0"
